$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily data rows (row 386 .. 464), continuing the date series from row 385
# (date serial 44459, 2021-09-20) through 2021-12-08 (serial 44538).
# "aggiornamento fino a 8/12" -> update through 8 December.
$startRow = 386
$startDate = 44460

$newData = @(
    @(0, 5, 32.29348317509527),
    @(0, 5, 32.29348317509527),
    @(0, 5, 32.29348317509527),
    @(0, 4, 25.83478654007622),
    @(1, 5, 32.29348317509527),
    @(0, 4, 25.83478654007622),
    @(1, 2, 12.91739327003811),
    @(0, 2, 12.91739327003811),
    @(0, 2, 12.91739327003811),
    @(0, 2, 12.91739327003811),
    @(0, 2, 12.91739327003811),
    @(1, 2, 12.91739327003811),
    @(0, 2, 12.91739327003811),
    @(0, 1, 6.458696635019054),
    @(0, 1, 6.458696635019054),
    @(0, 1, 6.458696635019054),
    @(1, 2, 12.91739327003811),
    @(0, 2, 12.91739327003811),
    @(0, 1, 6.458696635019054),
    @(0, 1, 6.458696635019054),
    @(0, 1, 6.458696635019054),
    @(0, 1, 6.458696635019054),
    @(0, 1, 6.458696635019054),
    @(0, 0, 0),
    @(0, 0, 0),
    @(0, 0, 0),
    @(0, 0, 0),
    @(0, 0, 0),
    @(0, 0, 0),
    @(0, 0, 0),
    @(0, 0, 0),
    @(1, 1, 6.458696635019054),
    @(0, 1, 6.458696635019054),
    @(0, 1, 6.458696635019054),
    @(2, 3, 19.37608990505716),
    @(0, 3, 19.37608990505716),
    @(0, 3, 19.37608990505716),
    @(1, 4, 25.83478654007622),
    @(1, 4, 25.83478654007622),
    @(1, 5, 32.29348317509527),
    @(4, 9, 58.12826971517148),
    @(2, 9, 58.12826971517148),
    @(0, 9, 58.12826971517148),
    @(0, 9, 58.12826971517148),
    @(0, 8, 51.66957308015243),
    @(0, 7, 45.21087644513337),
    @(4, 10, 64.58696635019054),
    @(0, 6, 38.75217981011431),
    @(1, 5, 32.29348317509527),
    @(3, 8, 51.66957308015243),
    @(0, 8, 51.66957308015243),
    @(1, 9, 58.12826971517148),
    @(4, 13, 83.96305625524769),
    @(0, 9, 58.12826971517148),
    @(0, 9, 58.12826971517148),
    @(0, 8, 51.66957308015243),
    @(15, 20, 129.1739327003811),
    @(0, 20, 129.1739327003811),
    @(0, 19, 122.715236065362),
    @(0, 15, 96.8804495252858),
    @(0, 15, 96.8804495252858),
    @(0, 15, 96.8804495252858),
    @(5, 20, 129.1739327003811),
    @(1, 6, 38.75217981011431),
    @(18, 24, 155.0087192404573),
    @(5, 29, 187.3022024155526),
    @(4, 33, 213.1369889556288),
    @(1, 34, 219.5956855906478),
    @(3, 37, 238.9717754957049),
    @(4, 36, 232.5130788606859),
    @(0, 35, 226.0543822256668),
    @(0, 17, 109.7978427953239),
    @(4, 16, 103.3391461603049),
    @(1, 13, 83.96305625524769),
    @(3, 15, 96.8804495252858),
    @(1, 13, 83.96305625524769),
    @(11, 20, 129.1739327003811),
    @(0, 20, 129.1739327003811),
    @(2, 22, 142.0913259704192)
)

for ($i = 0; $i -lt $newData.Count; $i++) {
    $r = $startRow + $i
    $entry = $newData[$i]

    # Column A: date, formatted/bordered like the existing date column (copy style from the row above)
    $ws.Range("A" + ($r - 1)).Copy($ws.Range("A" + $r))
    $ws.Cells.Item($r, 1).Value = $startDate + $i

    $ws.Cells.Item($r, 2).Value = $entry[0]
    $ws.Cells.Item($r, 3).Value = $entry[1]
    $ws.Cells.Item($r, 4).Value = $entry[2]
}

Write-Host "Inserted $($newData.Count) rows, last row: $($startRow + $newData.Count - 1)"
